# SpecFlowExample.feature.xlsx
#
# "Create CommonOutputSteps and reuse through scenarios"
#
# The two scenario sheets (MultiplicationTests / AdditioTests) each had
# their own bespoke "Then" step ("Output Multiplication is <MultiplicationOut>"
# / "Output Addition is <AdditionOut>"). The commit factors that out into a
# single shared/common output step ("Result is <ResultOut>") reused by both
# sheets, which is why the same three strings now appear identically on
# both worksheets (B3 / C3 / D7).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # MultiplicationTests
$ws2 = $wb.Worksheets.Item(2)   # AdditioTests

# --- Reposition the saved window (best effort; some hosts don't expose this) ---
$win = $wb.Windows.Item(1)
$win.Left = 2040
$win.Top  = 108

# --- Replace the per-sheet output step text with the shared "Result" step ---

# MultiplicationTests: "Output Multiplication is" / "<MultiplicationOut>" / "MultiplicationOut"
$ws1.Range("B3").Value = "Result is"
$ws1.Range("C3").Value = "<ResultOut>"
$ws1.Range("D7").Value = "ResultOut"

# AdditioTests: "Output Addition is" / "<AdditionOut>" / "AdditionOut"
$ws2.Range("B3").Value = "Result is"
$ws2.Range("C3").Value = "<ResultOut>"
$ws2.Range("D7").Value = "ResultOut"

# --- Move the saved cell selection on each sheet from D11 to B4 ---
$ws1.Range("B4").Select() | Out-Null
$ws2.Activate() | Out-Null
$ws2.Range("B4").Select() | Out-Null
